$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.398.41"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -2.97%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.462.57"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -2.98%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.64"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "93.95"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -6.67%  "
$ws.Range("E7").Value = "  -3.47%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -4.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.36"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -7.26%  "
$ws.Range("E11").Value = "  -2.94%  "
$ws.Range("E12").Value = "  -0.11%  "
$ws.Range("E13").Value = "  -5.76%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.843.66"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -3.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.474.37"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.67%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.49"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -9.18%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.788"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -3.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "41.403.06"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -2.97%  "
$ws.Range("E19").Value = "  -6.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0914"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -4.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.50"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -6.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.70"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "238.13"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.44%  "
$ws.Range("E24").Value = "  -4.02%  "
$ws.Range("E25").Value = "  -5.91%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.75"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -5.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.70"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -4.58%  "
$ws.Range("E30").Value = "  -7.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "153.02"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.38%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.62"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.90%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.63"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.55"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -6.65%  "
$ws.Range("E35").Value = "  -5.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.01"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -4.87%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.16"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -6.38%  "
$ws.Range("E38").Value = "  -7.38%  "
$ws.Range("E39").Value = "  -6.50%  "
$ws.Range("E40").Value = "  -3.99%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.11"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -5.27%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.26"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -3.11%  "
$ws.Range("E43").Value = "  +0.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.976.52"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0284"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -4.89%  "
$ws.Range("E46").Value = "  -7.80%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.77"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.17%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "76.89"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -5.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "97.32"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -3.75%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "69.12"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -4.64%  "
$ws.Range("E51").Value = "  -6.61%  "
